# Generate Report for Handoff
#
# A new handoff run produced a new GUID-named source file and new xliff
# hashes/timestamps. Update the Overview / zh-cn / de-de sheets (and the
# matching hyperlink display text) to reflect the new handoff artifacts.

$wb = $excel.ActiveWorkbook

$oldGuid = "74e31aec-6fb8-412e-b720-0dea785885ca"
$newGuid = "5d5eae22-f549-4ce7-bf34-414228d83089"

$newXlfHash = "2faf29f3b6f16e2861c5e7670cc7c6cb6ec7066e"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61bf03e54e115cf8467c7a0506452a5db31e8e14/e2e/"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-05 15:12:00"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $baseUrl + "$oldGuid.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-05 15:11:55"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $baseUrl + "$oldGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md")

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $baseUrl + "$oldGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md")
